$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Text)
    $escaped = $Text.Replace('"', '""')
    $ws.Range("ZZ1").Formula = '="' + $escaped + '"'
    $ws.Range("ZZ1").Copy() | Out-Null
    $ws.Range($CellRef).PasteSpecial(-4163) | Out-Null
    $ws.Range("ZZ1").ClearContents()
}

$ws.Range("D2").Value = "63.750.69"
$ws.Range("E2").Value = "  -6.35%  "
$ws.Range("D3").Value = "3.299.05"
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "183.04"
$ws.Range("E5").Value = "  -11.04%  "
Set-TextValue "D6" "521.11"
$ws.Range("E6").Value = "  -6.63%  "
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("D8").Value = "3.293.82"
$ws.Range("E8").Value = "  -7.54%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("E10").Value = "  -7.25%  "
Set-TextValue "D11" "58.77"
$ws.Range("E11").Value = "  -6.32%  "
$ws.Range("E12").Value = "  -8.72%  "
$ws.Range("E13").Value = "  -7.35%  "
Set-TextValue "D14" "9.16"
$ws.Range("E14").Value = "  -8.72%  "
$ws.Range("D15").Value = "3.812.27"
$ws.Range("E15").Value = "  -7.75%  "
$ws.Range("E16").Value = "  -5.38%  "
$ws.Range("D17").Value = "3.290.06"
$ws.Range("E17").Value = "  -7.58%  "
$ws.Range("E18").Value = "  -6.32%  "
$ws.Range("D19").Value = "63.640.92"
$ws.Range("E19").Value = "  -6.12%  "
$ws.Range("E20").Value = "  -8.99%  "
Set-TextValue "D21" "0.953"
$ws.Range("E21").Value = "  -9.43%  "
Set-TextValue "D22" "372.90"
$ws.Range("E22").Value = "  -5.49%  "
Set-TextValue "D23" "11.30"
$ws.Range("E23").Value = "  -7.36%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D24" "3.71"
$ws.Range("E24").Value = "  -9.15%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D25" "80.32"
$ws.Range("E25").Value = "  -4.21%  "
Set-TextValue "D26" "3.84"
$ws.Range("E26").Value = "  -0.11%  "
Set-TextValue "D27" "5.98"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("E28").Value = "  -6.81%  "
$ws.Range("E29").Value = "  -7.22%  "
$ws.Range("E30").Value = "  -7.67%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D31" "651.57"
$ws.Range("E31").Value = "  -9.19%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D32" "28.66"
$ws.Range("E32").Value = "  -8.07%  "
$ws.Range("E33").Value = "  -10.22%  "
$ws.Range("E34").Value = "  -6.30%  "
Set-TextValue "D35" "59.41"
$ws.Range("E35").Value = "  -6.71%  "
$ws.Range("E36").Value = "  -5.41%  "
$ws.Range("E37").Value = "  -0.05%  "
Set-TextValue "D38" "0.392"
$ws.Range("E38").Value = "  -5.62%  "
Set-TextValue "D39" "36.28"
$ws.Range("E39").Value = "  -11.07%  "
Set-TextValue "D40" "0.997"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "2.993.17"
$ws.Range("E41").Value = "  -4.57%  "
Set-TextValue "D42" "0.126"
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("D43").Value = "0.0₃0654"
$ws.Range("E43").Value = "  -10.01%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D44" "2.44"
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D45" "2.69"
$ws.Range("E45").Value = "  -15.90%  "
Set-TextValue "D46" "0.0391"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D47" "2.82"
$ws.Range("E47").Value = "  +5.84%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D48" "2.58"
$ws.Range("E48").Value = "  -5.73%  "
Set-TextValue "D49" "0.125"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("E50").Value = "  -3.48%  "
Set-TextValue "D51" "2.50"
$ws.Range("E51").Value = "  -18.71%  "
